$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 66.47695399999999
$ws.Range("H2").Value = 199.430862
$ws.Range("I2").Value = 0.04311983106164722
$ws.Range("J2").Value = 0.04311983106164721
$ws.Range("M2").Value = 14.89002333333333
$ws.Range("N2").Value = 44.67007
$ws.Range("O2").Value = 0.1194491234330596
$ws.Range("P2").Value = 0.1194491234330597
$ws.Range("Q2").Value = 989.8433961889266
$ws.Range("R2").Value = 8908.59056570034
$ws.Range("S2").Value = 0.005150626022895378
$ws.Range("T2").Value = 0.005150626022895378
$ws.Range("G3").Value = 66.47695399999999
$ws.Range("H3").Value = 199.430862
$ws.Range("I3").Value = 0.04311983106164722
$ws.Range("J3").Value = 0.04311983106164721
$ws.Range("O3").Value = 0.05148509068166413
$ws.Range("P3").Value = 0.05148509068166414
$ws.Range("Q3").Value = 426.6433737539559
$ws.Range("R3").Value = 3839.790363785603
$ws.Range("S3").Value = 0.002220028412386945
$ws.Range("T3").Value = 0.002220028412386945
$ws.Range("G4").Value = 66.47695399999999
$ws.Range("H4").Value = 199.430862
$ws.Range("I4").Value = 0.04311983106164722
$ws.Range("J4").Value = 0.04311983106164721
$ws.Range("M4").Value = 66.36284166666667
$ws.Range("N4").Value = 199.088525
$ws.Range("O4").Value = 0.5323687604884161
$ws.Range("P4").Value = 0.5323687604884162
$ws.Range("Q4").Value = 4411.599572784283
$ws.Range("R4").Value = 39704.39615505855
$ws.Range("S4").Value = 0.02295565101475903
$ws.Range("T4").Value = 0.02295565101475903
$ws.Range("G5").Value = 66.47695399999999
$ws.Range("H5").Value = 199.430862
$ws.Range("I5").Value = 0.04311983106164722
$ws.Range("J5").Value = 0.04311983106164721
$ws.Range("M5").Value = 3.521285666666667
$ws.Range("N5").Value = 10.563857
$ws.Range("O5").Value = 0.02824807435318976
$ws.Range("P5").Value = 0.02824807435318976
$ws.Range("Q5").Value = 234.0843452838593
$ws.Range("R5").Value = 2106.759107554734
$ws.Range("S5").Value = 0.001218052193926392
$ws.Range("T5").Value = 0.001218052193926392
$ws.Range("G6").Value = 66.47695399999999
$ws.Range("H6").Value = 199.430862
$ws.Range("I6").Value = 0.04311983106164722
$ws.Range("J6").Value = 0.04311983106164721
$ws.Range("M6").Value = 33.46371266666667
$ws.Range("N6").Value = 100.391138
$ws.Range("O6").Value = 0.2684489510436703
$ws.Range("P6").Value = 0.2684489510436703
$ws.Range("Q6").Value = 2224.565687611217
$ws.Range("R6").Value = 20021.09118850096
$ws.Range("S6").Value = 0.01157547341767947
$ws.Range("T6").Value = 0.01157547341767946
$ws.Range("I7").Value = 0.8830494168872806
$ws.Range("J7").Value = 0.8830494168872804
$ws.Range("M7").Value = 14.89002333333333
$ws.Range("N7").Value = 44.67007
$ws.Range("O7").Value = 0.1194491234330596
$ws.Range("P7").Value = 0.1194491234330597
$ws.Range("Q7").Value = 20270.96610292161
$ws.Range("R7").Value = 182438.6949262945
$ws.Range("S7").Value = 0.1054794787952601
$ws.Range("T7").Value = 0.1054794787952601
$ws.Range("I8").Value = 0.8830494168872806
$ws.Range("J8").Value = 0.8830494168872804
$ws.Range("O8").Value = 0.05148509068166413
$ws.Range("P8").Value = 0.05148509068166414
$ws.Range("S8").Value = 0.04546387930483228
$ws.Range("T8").Value = 0.04546387930483228
$ws.Range("I9").Value = 0.8830494168872806
$ws.Range("J9").Value = 0.8830494168872804
$ws.Range("M9").Value = 66.36284166666667
$ws.Range("N9").Value = 199.088525
$ws.Range("O9").Value = 0.5323687604884161
$ws.Range("P9").Value = 0.5323687604884162
$ws.Range("Q9").Value = 90344.98360436107
$ws.Range("R9").Value = 813104.8524392496
$ws.Range("S9").Value = 0.4701079235183002
$ws.Range("T9").Value = 0.4701079235183002
$ws.Range("I10").Value = 0.8830494168872806
$ws.Range("J10").Value = 0.8830494168872804
$ws.Range("M10").Value = 3.521285666666667
$ws.Range("N10").Value = 10.563857
$ws.Range("O10").Value = 0.02824807435318976
$ws.Range("P10").Value = 0.02824807435318976
$ws.Range("Q10").Value = 4793.804602569711
$ws.Range("R10").Value = 43144.2414231274
$ws.Range("S10").Value = 0.02494444558577276
$ws.Range("T10").Value = 0.02494444558577276
$ws.Range("I11").Value = 0.8830494168872806
$ws.Range("J11").Value = 0.8830494168872804
$ws.Range("M11").Value = 33.46371266666667
$ws.Range("N11").Value = 100.391138
$ws.Range("O11").Value = 0.2684489510436703
$ws.Range("P11").Value = 0.2684489510436703
$ws.Range("Q11").Value = 45556.79799543017
$ws.Range("R11").Value = 410011.1819588715
$ws.Range("S11").Value = 0.2370536896831152
$ws.Range("T11").Value = 0.2370536896831151
$ws.Range("G12").Value = 44.831112
$ws.Range("H12").Value = 134.493336
$ws.Range("I12").Value = 0.02907940059566787
$ws.Range("J12").Value = 0.02907940059566786
$ws.Range("M12").Value = 14.89002333333333
$ws.Range("N12").Value = 44.67007
$ws.Range("O12").Value = 0.1194491234330596
$ws.Range("P12").Value = 0.1194491234330597
$ws.Range("Q12").Value = 667.53630373928
$ws.Range("R12").Value = 6007.82673365352
$ws.Range("S12").Value = 0.00347350891111132
$ws.Range("T12").Value = 0.003473508911111319
$ws.Range("G13").Value = 44.831112
$ws.Range("H13").Value = 134.493336
$ws.Range("I13").Value = 0.02907940059566787
$ws.Range("J13").Value = 0.02907940059566786
$ws.Range("O13").Value = 0.05148509068166413
$ws.Range("P13").Value = 0.05148509068166414
$ws.Range("Q13").Value = 287.722221340368
$ws.Range("R13").Value = 2589.499992063312
$ws.Range("S13").Value = 0.001497155576636398
$ws.Range("T13").Value = 0.001497155576636398
$ws.Range("G14").Value = 44.831112
$ws.Range("H14").Value = 134.493336
$ws.Range("I14").Value = 0.02907940059566787
$ws.Range("J14").Value = 0.02907940059566786
$ws.Range("M14").Value = 66.36284166666667
$ws.Range("N14").Value = 199.088525
$ws.Range("O14").Value = 0.5323687604884161
$ws.Range("P14").Value = 0.5323687604884162
$ws.Range("Q14").Value = 2975.1199873966
$ws.Range("R14").Value = 26776.0798865694
$ws.Range("S14").Value = 0.01548096445086181
$ws.Range("T14").Value = 0.01548096445086181
$ws.Range("G15").Value = 44.831112
$ws.Range("H15").Value = 134.493336
$ws.Range("I15").Value = 0.02907940059566787
$ws.Range("J15").Value = 0.02907940059566786
$ws.Range("M15").Value = 3.521285666666667
$ws.Range("N15").Value = 10.563857
$ws.Range("O15").Value = 0.02824807435318976
$ws.Range("P15").Value = 0.02824807435318976
$ws.Range("Q15").Value = 157.863152106328
$ws.Range("R15").Value = 1420.768368956952
$ws.Range("S15").Value = 0.0008214370701726166
$ws.Range("T15").Value = 0.0008214370701726165
$ws.Range("G16").Value = 44.831112
$ws.Range("H16").Value = 134.493336
$ws.Range("I16").Value = 0.02907940059566787
$ws.Range("J16").Value = 0.02907940059566786
$ws.Range("M16").Value = 33.46371266666667
$ws.Range("N16").Value = 100.391138
$ws.Range("O16").Value = 0.2684489510436703
$ws.Range("P16").Value = 0.2684489510436703
$ws.Range("Q16").Value = 1500.215450495152
$ws.Range("R16").Value = 13501.93905445637
$ws.Range("S16").Value = 0.00780633458688572
$ws.Range("T16").Value = 0.007806334586885719
$ws.Range("G17").Value = 52.83062100000001
$ws.Range("H17").Value = 158.491863
$ws.Range("I17").Value = 0.0342682285413064
$ws.Range("J17").Value = 0.03426822854130639
$ws.Range("M17").Value = 14.89002333333333
$ws.Range("N17").Value = 44.67007
$ws.Range("O17").Value = 0.1194491234330596
$ws.Range("P17").Value = 0.1194491234330597
$ws.Range("Q17").Value = 786.6491794044902
$ws.Range("R17").Value = 7079.842614640412
$ws.Range("S17").Value = 0.004093309860862805
$ws.Range("T17").Value = 0.004093309860862805
$ws.Range("G18").Value = 52.83062100000001
$ws.Range("H18").Value = 158.491863
$ws.Range("I18").Value = 0.0342682285413064
$ws.Range("J18").Value = 0.03426822854130639
$ws.Range("O18").Value = 0.05148509068166413
$ws.Range("P18").Value = 0.05148509068166414
$ws.Range("Q18").Value = 339.062382144594
$ws.Range("R18").Value = 3051.561439301346
$ws.Range("S18").Value = 0.001764302853949151
$ws.Range("T18").Value = 0.001764302853949151
$ws.Range("G19").Value = 52.83062100000001
$ws.Range("H19").Value = 158.491863
$ws.Range("I19").Value = 0.0342682285413064
$ws.Range("J19").Value = 0.03426822854130639
$ws.Range("M19").Value = 66.36284166666667
$ws.Range("N19").Value = 199.088525
$ws.Range("O19").Value = 0.5323687604884161
$ws.Range("P19").Value = 0.5323687604884162
$ws.Range("Q19").Value = 3505.990136574675
$ws.Range("R19").Value = 31553.91122917208
$ws.Range("S19").Value = 0.01824333435266905
$ws.Range("T19").Value = 0.01824333435266905
$ws.Range("G20").Value = 52.83062100000001
$ws.Range("H20").Value = 158.491863
$ws.Range("I20").Value = 0.0342682285413064
$ws.Range("J20").Value = 0.03426822854130639
$ws.Range("M20").Value = 3.521285666666667
$ws.Range("N20").Value = 10.563857
$ws.Range("O20").Value = 0.02824807435318976
$ws.Range("P20").Value = 0.02824807435318976
$ws.Range("Q20").Value = 186.031708488399
$ws.Range("R20").Value = 1674.285376395591
$ws.Range("S20").Value = 0.0009680114677869227
$ws.Range("T20").Value = 0.0009680114677869224
$ws.Range("G21").Value = 52.83062100000001
$ws.Range("H21").Value = 158.491863
$ws.Range("I21").Value = 0.0342682285413064
$ws.Range("J21").Value = 0.03426822854130639
$ws.Range("M21").Value = 33.46371266666667
$ws.Range("N21").Value = 100.391138
$ws.Range("O21").Value = 0.2684489510436703
$ws.Range("P21").Value = 0.2684489510436703
$ws.Range("Q21").Value = 1767.908721145566
$ws.Range("R21").Value = 15911.1784903101
$ws.Range("S21").Value = 0.009199270006038467
$ws.Range("T21").Value = 0.009199270006038464
$ws.Range("G22").Value = 16.16161433333333
$ws.Range("H22").Value = 48.484843
$ws.Range("I22").Value = 0.01048312291409786
$ws.Range("J22").Value = 0.01048312291409786
$ws.Range("M22").Value = 14.89002333333333
$ws.Range("N22").Value = 44.67007
$ws.Range("O22").Value = 0.1194491234330596
$ws.Range("P22").Value = 0.1194491234330597
$ws.Range("Q22").Value = 240.6468145276678
$ws.Range("R22").Value = 2165.82133074901
$ws.Range("S22").Value = 0.001252199842930012
$ws.Range("T22").Value = 0.001252199842930012
$ws.Range("G23").Value = 16.16161433333333
$ws.Range("H23").Value = 48.484843
$ws.Range("I23").Value = 0.01048312291409786
$ws.Range("J23").Value = 0.01048312291409786
$ws.Range("O23").Value = 0.05148509068166413
$ws.Range("P23").Value = 0.05148509068166414
$ws.Range("Q23").Value = 103.7238508925006
$ws.Range("R23").Value = 933.514658032506
$ws.Range("S23").Value = 0.0005397245338593597
$ws.Range("T23").Value = 0.0005397245338593597
$ws.Range("G24").Value = 16.16161433333333
$ws.Range("H24").Value = 48.484843
$ws.Range("I24").Value = 0.01048312291409786
$ws.Range("J24").Value = 0.01048312291409786
$ws.Range("M24").Value = 66.36284166666667
$ws.Range("N24").Value = 199.088525
$ws.Range("O24").Value = 0.5323687604884161
$ws.Range("P24").Value = 0.5323687604884162
$ws.Range("Q24").Value = 1072.530653080731
$ws.Range("R24").Value = 9652.775877726575
$ws.Range("S24").Value = 0.005580887151825992
$ws.Range("T24").Value = 0.005580887151825992
$ws.Range("G25").Value = 16.16161433333333
$ws.Range("H25").Value = 48.484843
$ws.Range("I25").Value = 0.01048312291409786
$ws.Range("J25").Value = 0.01048312291409786
$ws.Range("M25").Value = 3.521285666666667
$ws.Range("N25").Value = 10.563857
$ws.Range("O25").Value = 0.02824807435318976
$ws.Range("P25").Value = 0.02824807435318976
$ws.Range("Q25").Value = 56.90966090216123
$ws.Range("R25").Value = 512.186948119451
$ws.Range("S25").Value = 0.0002961280355310638
$ws.Range("T25").Value = 0.0002961280355310637
$ws.Range("G26").Value = 16.16161433333333
$ws.Range("H26").Value = 48.484843
$ws.Range("I26").Value = 0.01048312291409786
$ws.Range("J26").Value = 0.01048312291409786
$ws.Range("M26").Value = 33.46371266666667
$ws.Range("N26").Value = 100.391138
$ws.Range("O26").Value = 0.2684489510436703
$ws.Range("P26").Value = 0.2684489510436703
$ws.Range("Q26").Value = 540.8276182801482
$ws.Range("R26").Value = 4867.448564521334
$ws.Range("S26").Value = 0.002814183349951436
$ws.Range("T26").Value = 0.002814183349951435
